$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -1

$ws.Range("B3").Value = -25
$ws.Range("C3").Value = -9

$ws.Range("C4").Value = -14

$ws.Range("C5").Value = 9

$ws.Range("C6").Value = 43

$ws.Range("C7").Value = 153

$ws.Range("C8").Value = 275

$ws.Range("C9").Value = 244

$ws.Range("C11").Value = 140

$ws.Range("C18").Value = -29
